# Applies the Serie A 2023-2024 update script (run 26-11-2023 20:30):
#  - Swaps the two Matchday-7 fixtures that were logged out of kickoff-time order
#    (Inter-Sassuolo / Lazio-Torino), rows 57/58.
#  - Swaps the two Matchday-12 fixtures that were logged out of kickoff-time order
#    (Udinese-Atalanta / Fiorentina-Bologna), rows 118/119.
#  - Appends 7 newly played fixtures (rows 122-128).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-RowData($rowA, $rowB) {
    $valsA = @{}
    $valsB = @{}
    foreach ($c in $dataCols) {
        $valsA[$c] = $ws.Range("$c$rowA").Value()
        $valsB[$c] = $ws.Range("$c$rowB").Value()
    }
    foreach ($c in $dataCols) {
        $ws.Range("$c$rowA").Value = $valsB[$c]
        $ws.Range("$c$rowB").Value = $valsA[$c]
    }
}

# --- Swap the two mis-ordered fixtures around row 57/58 ---
Swap-RowData 57 58

# --- Swap the two mis-ordered fixtures around row 118/119 ---
Swap-RowData 118 119

# --- Append the 7 newly played fixtures as rows 122-128 ---
function Set-IndexCell($row, $indice) {
    $cell = $ws.Range("A$row")
    $cell.Value = $indice
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
}

function Set-DateCell($row, $serial) {
    $cell = $ws.Range("E$row")
    $cell.Value = $serial
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

function Set-MatchRow($row, $indice, $home, $homeGoals, $away, $awayGoals,
                       $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
                       $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
                       $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt,
                       $matchDateSerial, $url) {
    Set-IndexCell $row $indice
    $ws.Range("B$row").Value = "italy"
    $ws.Range("C$row").Value = "serie-a"
    $ws.Range("D$row").Value = "2023-2024"
    Set-DateCell $row $matchDateSerial
    $ws.Range("F$row").Value = $home
    $ws.Range("G$row").Value = $homeGoals
    $ws.Range("H$row").Value = $away
    $ws.Range("I$row").Value = $awayGoals
    $ws.Range("J$row").Value = $homeOpenOdds
    $ws.Range("K$row").Value = $homeOpenDt
    $ws.Range("L$row").Value = $homeCloseOdds
    $ws.Range("M$row").Value = $homeCloseDt
    $ws.Range("N$row").Value = $drawOpenOdds
    $ws.Range("O$row").Value = $drawOpenDt
    $ws.Range("P$row").Value = $drawCloseOdds
    $ws.Range("Q$row").Value = $drawCloseDt
    $ws.Range("R$row").Value = $awayOpenOdds
    $ws.Range("S$row").Value = $awayOpenDt
    $ws.Range("T$row").Value = $awayCloseOdds
    $ws.Range("U$row").Value = $awayCloseDt
    $ws.Range("V$row").Value = $url
}

Set-MatchRow 122 121 "Salernitana" 2 "Lazio" 1 `
    3.75 "05/11/2023 11:03" 4.58 "25/11/2023 14:59" `
    3.47 "05/11/2023 11:03" 3.46 "25/11/2023 14:59" `
    2.06 "05/11/2023 11:03" 1.92 "25/11/2023 14:57" `
    45255.625 "https://www.betexplorer.com/football/italy/serie-a/salernitana-lazio/8x6P9U9q/"

Set-MatchRow 123 122 "Atalanta" 1 "Napoli" 2 `
    2.6 "05/11/2023 11:03" 2.53 "25/11/2023 17:59" `
    3.41 "05/11/2023 11:03" 3.42 "25/11/2023 17:57" `
    2.77 "05/11/2023 11:03" 2.94 "25/11/2023 17:59" `
    45255.75 "https://www.betexplorer.com/football/italy/serie-a/atalanta-napoli/UVIceVfj/"

Set-MatchRow 124 123 "AC Milan" 1 "Fiorentina" 0 `
    1.72 "05/11/2023 11:03" 2.13 "25/11/2023 20:39" `
    3.84 "05/11/2023 11:03" 3.37 "25/11/2023 20:35" `
    5.02 "05/11/2023 11:03" 3.83 "25/11/2023 20:43" `
    45255.86458333334 "https://www.betexplorer.com/football/italy/serie-a/ac-milan-fiorentina/xfFgFjnM/"

Set-MatchRow 125 124 "Cagliari" 1 "Monza" 1 `
    2.68 "05/11/2023 11:03" 2.72 "26/11/2023 12:29" `
    3.3 "05/11/2023 11:03" 3.25 "26/11/2023 12:21" `
    2.76 "05/11/2023 11:03" 2.84 "26/11/2023 12:29" `
    45256.52083333334 "https://www.betexplorer.com/football/italy/serie-a/cagliari-monza/8ENxJYnc/"

Set-MatchRow 126 125 "Empoli" 3 "Sassuolo" 4 `
    2.81 "05/11/2023 11:03" 2.67 "26/11/2023 14:59" `
    3.52 "05/11/2023 11:03" 3.55 "26/11/2023 14:59" `
    2.51 "05/11/2023 11:03" 2.69 "26/11/2023 14:59" `
    45256.625 "https://www.betexplorer.com/football/italy/serie-a/empoli-sassuolo/n5MtIh23/"

Set-MatchRow 127 126 "Frosinone" 2 "Genoa" 1 `
    2.59 "05/11/2023 11:03" 2.48 "26/11/2023 14:59" `
    3.24 "05/11/2023 11:03" 3.22 "26/11/2023 14:57" `
    2.91 "05/11/2023 11:03" 3.18 "26/11/2023 14:59" `
    45256.625 "https://www.betexplorer.com/football/italy/serie-a/frosinone-genoa/tSLpHCH9/"

Set-MatchRow 128 127 "AS Roma" 3 "Udinese" 1 `
    1.66 "05/11/2023 11:03" 1.66 "26/11/2023 17:32" `
    3.83 "05/11/2023 11:03" 3.83 "26/11/2023 17:32" `
    5.57 "05/11/2023 11:03" 6 "26/11/2023 17:59" `
    45256.75 "https://www.betexplorer.com/football/italy/serie-a/as-roma-udinese/GpEcEA1S/"

Write-Host "Edit complete"
